$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.41"
$ws.Range("E2").Value = "'-0.64%"
$ws.Range("D3").Value = "'28.30"
$ws.Range("E3").Value = "'-4.59%"
$ws.Range("D4").Value = "'5.254"
$ws.Range("E4").Value = "'1.68%"
$ws.Range("D5").Value = "'0.05699"
$ws.Range("E5").Value = "'-0.55%"
$ws.Range("D6").Value = "'6.632"
$ws.Range("E6").Value = "'0.87%"
$ws.Range("D7").Value = "'3.204"
$ws.Range("E7").Value = "'3.42%"
$ws.Range("D8").Value = "'0.8510"
$ws.Range("E8").Value = "'-0.71%"
$ws.Range("D9").Value = "'0.9077"
$ws.Range("E9").Value = "'4.60%"
$ws.Range("D10").Value = "'0.1368"
$ws.Range("E10").Value = "'0.28%"
$ws.Range("D11").Value = "'0.07078"
$ws.Range("E11").Value = "'-0.08%"
$ws.Range("D12").Value = "'0.03186"
$ws.Range("E12").Value = "'8.75%"
$ws.Range("D13").Value = "'0.09225"
$ws.Range("E13").Value = "'-1.66%"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'1.03%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005947"
$ws.Range("E15").Value = "'-1.14%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005934"
$ws.Range("E16").Value = "'-2.44%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.491"
$ws.Range("E17").Value = "'0.04%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.188"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("E19").Value = "'-0.46%"
$ws.Range("D20").Value = "'0.03314"
$ws.Range("E20").Value = "'-2.41%"
$ws.Range("E21").Value = "'-2.05%"
$ws.Range("D22").Value = "'3.522"
$ws.Range("E22").Value = "'1.73%"
$ws.Range("D23").Value = "'0.04073"
$ws.Range("E23").Value = "'-1.37%"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-0.34%"
$ws.Range("D26").Value = "'0.004152"
$ws.Range("E26").Value = "'-17.10%"
$ws.Range("E27").Value = "'-0.81%"
$ws.Range("D40").Value = "'0.03824"
$ws.Range("E40").Value = "'1.80%"
$ws.Range("D41").Value = "'0.1067"
$ws.Range("E41").Value = "'-0.60%"
$ws.Range("D42").Value = "'0.003737"
$ws.Range("E42").Value = "'-34.85%"
$ws.Range("D43").Value = "'0.002489"
$ws.Range("E43").Value = "'1.65%"
$ws.Range("D44").Value = "'0.009161"
$ws.Range("E44").Value = "'7.82%"
$ws.Range("D45").Value = "'0.00005269"
$ws.Range("E45").Value = "'0.47%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("D47").Value = "'0.1050"
$ws.Range("E47").Value = "'62.33%"
$ws.Range("D48").Value = "'0.002269"
$ws.Range("E48").Value = "'0.48%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.02%"
